$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175815343856812
$ws.Range("B1").Value = 2.216834306716919
$ws.Range("C1").Value = 4.500626564025879
$ws.Range("D1").Value = 2.660823106765747
$ws.Range("E1").Value = 1.224990487098694
